# Applies the "feat: add 2022-Q3 data" change:
#  1. Insert a new worksheet "2022-Q3" right after "总计", pushing all the
#     other quarter sheets (2022-Q2 .. 2021-Q1) down by one tab position.
#  2. Populate the new "2022-Q3" sheet with the fund-holdings table,
#     copying the header/index-column formatting from the existing
#     "2022-Q2" sheet (which has the exact same table layout).
#  3. Insert a new summary row for "2022-Q3" at the top of the "总计"
#     sheet's data (row 2), shifting the existing rows down and
#     renumbering the index column, without disturbing existing styles.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after "总计"
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2. Populate "2022-Q3" - reuse the formatting from "2022-Q2" (same layout:
#    bold/centered header row, bordered index column A) before writing data.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")

$template.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$template.Range("A2").Copy()
$q3.Range("A2:A11").PasteSpecial(-4122)

# Header row
$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Data rows: index, code, name, size, totalPosition, positionRatio, marketValue, rank
$q3Data = @(
    @(0, "519674", "银河创新成长混合A",               "142.10", "94.51", "7.99", "11.3538", 4),
    @(1, "014143", "银河创新成长混合C",               "20.25",  "94.51", "7.99", "1.6180",  4),
    @(2, "217021", "招商优势企业混合",                 "2.91",   "79.74", "5.08", "0.1478",  7),
    @(3, "000522", "华润元大信息传媒科技混合",         "1.42",   "71.23", "4.77", "0.0677",  4),
    @(4, "013340", "创金合信芯片产业股票C",            "0.93",   "90.74", "6.65", "0.0618",  4),
    @(5, "013339", "创金合信芯片产业股票A",            "0.91",   "90.74", "6.65", "0.0605",  4),
    @(6, "002772", "光大保德信产业新动力灵活配置混合", "0.25",   "87.41", "6.12", "0.0153",  7),
    @(7, "004890", "中邮健康文娱灵活配置混合",         "0.41",   "86.15", "3.71", "0.0152",  8),
    @(8, "004931", "华润元大价值优选混合C",            "0.15",   "66.44", "3.73", "0.0056",  10),
    @(9, "004930", "华润元大价值优选混合A",            "0.12",   "66.44", "3.73", "0.0045",  10)
)

$r = 2
foreach ($row in $q3Data) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = "'" + $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = "'" + $row[3]
    $q3.Cells.Item($r, 5).Value = "'" + $row[4]
    $q3.Cells.Item($r, 6).Value = "'" + $row[5]
    $q3.Cells.Item($r, 7).Value = "'" + $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3. Update the "总计" sheet: insert the 2022-Q3 summary row and shift the
#    rest of the table down by one row (bottom-up, reading via Value2 since
#    the Value getter is unreliable in this host).
# ---------------------------------------------------------------------------

# Extend column A's index-cell styling down onto the newly used row 8
$zongji.Range("A7").Copy()
$zongji.Range("A8").PasteSpecial(-4122)

for ($r = 7; $r -ge 2; $r--) {
    $zongji.Cells.Item($r + 1, 2).Value = $zongji.Cells.Item($r, 2).Value2
    $zongji.Cells.Item($r + 1, 3).Value = $zongji.Cells.Item($r, 3).Value2
    $zongji.Cells.Item($r + 1, 4).Value = $zongji.Cells.Item($r, 4).Value2
}

$zongji.Cells.Item(2, 2).Value = "2022-Q3"
$zongji.Cells.Item(2, 3).Value = 10
$zongji.Cells.Item(2, 4).Value = 13.35

# Renumber the index column (0..6) for rows 2..8
for ($r = 2; $r -le 8; $r++) {
    $zongji.Cells.Item($r, 1).Value = $r - 2
}
